$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-11-11 Tuesday" "2025-11-12 Wednesday"

Replace-Text "64÷6=" "61÷7="
Replace-Text "54÷2=" "88÷6="
Replace-Text "76÷6=" "99÷4="
Replace-Text "71÷2=" "75÷7="
Replace-Text "24÷2=" "41÷8="
Replace-Text "77÷3=" "74÷7="
Replace-Text "61÷6=" "67÷6="
Replace-Text "39÷2=" "49÷4="
Replace-Text "90÷8=" "55÷5="
Replace-Text "78÷3=" "39÷9="
Replace-Text "87÷6=" "52÷5="
Replace-Text "66÷5=" "10÷6="
Replace-Text "48÷8=" "37÷4="
Replace-Text "33÷8=" "79÷8="
Replace-Text "94÷6=" "54÷9="
Replace-Text "73÷2=" "21÷4="
Replace-Text "51÷9=" "33÷4="
Replace-Text "23÷8=" "13÷4="
Replace-Text "45÷6=" "87÷2="
Replace-Text "54÷4=" "40÷3="
Replace-Text "81÷2=" "72÷2="
Replace-Text "24÷6=" "59÷2="
Replace-Text "21÷7=" "92÷8="
Replace-Text "93÷4=" "35÷6="
Replace-Text "46÷5=" "58÷2="
